$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the balance value for No.6 (row 6) from 500 to 428
$ws.Range("B6").Value = 428

# Update the selected cell to C16 to match the saved selection state
$ws.Range("C16").Select()
